# Edit the workbook per the commit: insert two new columns ("Cong phu phau 1"
# and "Cong phu phau 2") into sheet 1, add a Total row, bump two values in
# sheet 2, and add a brand-new "LUY KE NGAY" sheet with daily cumulative data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: DOANH SO CA NHAN
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("DOANH SỐ CÁ NHÂN")

# Insert a new column before current G ("Số lần phụ phẫu 2") to hold the new
# "Công phụ phẫu 1" values, and another new column before current H (now I)
# ("Doanh số thu nợ") to hold "Công phụ phẫu 2".
$ws1.Range("G1:G11").EntireColumn.Insert()
$ws1.Range("I1:I11").EntireColumn.Insert()

$ws1.Range("G1").Value = "Công phụ phẫu 1"
$ws1.Range("I1").Value = "Công phụ phẫu 2"

# Fill the new "Công phụ phẫu 1" column (G) with the per-employee values.
$ws1.Range("G2").Value = 0
$ws1.Range("G3").Value = 150000
$ws1.Range("G4").Value = 1850000
$ws1.Range("G5").Value = 0
$ws1.Range("G6").Value = 0
$ws1.Range("G7").Value = 0
$ws1.Range("G8").Value = 0
$ws1.Range("G9").Value = 1050000
$ws1.Range("G10").Value = 0
$ws1.Range("G11").Value = 50000

# Fill the new "Công phụ phẫu 2" column (I) - all zero in this period.
$ws1.Range("I2:I11").Value = 0

# Add the "Tổng" (total) row 12.
$ws1.Range("A12").Value = "Tổng"
$ws1.Range("B12").Value = 618100000
$ws1.Range("C12").Value = 0
$ws1.Range("D12").Value = 511100000
$ws1.Range("E12").Value = 0
$ws1.Range("F12").Value = 33
$ws1.Range("G12").Value = 3100000
$ws1.Range("H12").Value = 3
$ws1.Range("I12").Value = 0
$ws1.Range("J12").Value = 64200000

# ---------------------------------------------------------------------------
# Sheet 2: CHI TIEU
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("CHI TIÊU")
$ws2.Range("B10").Value = 16200000
$ws2.Range("B12").Value = 164473000

# ---------------------------------------------------------------------------
# Sheet 3 (new): LUY KE NGAY - inserted right after "CHI TIÊU"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "LŨY KẾ NGÀY"

$headers3 = @("Ngày", "Đơn giá", "Thanh toán lần đầu", "Số lượng đơn", "Thu nợ", "Lượng chi")
for ($c = 0; $c -lt $headers3.Length; $c++) {
    $ws3.Cells.Item(1, $c + 1).Value = $headers3[$c]
}

# Column A holds "MM-DD-YYYY" labels that must stay plain text (not be
# auto-converted to date serials) - format as Text before writing them.
$ws3.Range("A2:A28").NumberFormat = "@"

$data3 = @(
    @("06-03-2024", 73000000, 70000000, 8, 0, 1490000),
    @("06-04-2024", 20000000, 15000000, 1, 0, 0),
    @("06-05-2024", 19000000, 19000000, 3, 1000000, 6750000),
    @("06-06-2024", 46000000, 41000000, 2, 6500000, 1008000),
    @("06-07-2024", 4500000, 4500000, 1, 18500000, 7610000),
    @("06-08-2024", 26000000, 23000000, 2, 0, 0),
    @("06-09-2024", 72900000, 63900000, 5, 0, 3010000),
    @("06-10-2024", 79000000, 79000000, 4, 0, 4435000),
    @("06-11-2024", 30000000, 30000000, 1, 0, 11640000),
    @("06-13-2024", 1300000, 1300000, 1, 10000000, 5000000),
    @("06-14-2024", 35000000, 35000000, 2, 8000000, 16000000),
    @("06-15-2024", 9000000, 9000000, 2, 10200000, 1000000),
    @("06-16-2024", 29500000, 29500000, 4, 1000000, 5810000),
    @("06-17-2024", 83000000, 83000000, 3, 0, 20900000),
    @("06-18-2024", 7000000, 6500000, 3, 0, 4020000),
    @("06-19-2024", 18000000, 13000000, 3, 0, 490000),
    @("06-20-2024", 2000000, 2000000, 1, 0, 10100000),
    @("06-21-2024", 39400000, 39400000, 6, 1000000, 4346000),
    @("06-22-2024", 0, 0, 0, 1000000, 7500000),
    @("06-23-2024", 15000000, 10000000, 1, 0, 9400000),
    @("06-24-2024", 3000000, 3000000, 2, 0, 1350000),
    @("06-25-2024", 3500000, 3500000, 2, 0, 840000),
    @("06-26-2024", 2000000, 2000000, 1, 1000000, 11800000),
    @("06-27-2024", 0, 0, 0, 3000000, 0),
    @("06-28-2024", 0, 0, 0, 0, 7350000),
    @("06-29-2024", 0, 0, 0, 3000000, 0),
    @("06-30-2024", 0, 0, 0, 5000000, 22624000),
    @("Tổng", 618100000, 582600000, 58, 69200000, 164473000)
)

$r = 2
foreach ($row in $data3) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws3.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}
